$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("N10").Value = 64.42
$ws.Range("O10").Value = 4260
$ws.Range("N12").Value = 4267.99
$ws.Range("O12").Value = 62.9
